$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19 (shifts existing rows 19-28 down to 20-29),
# inheriting the row-19 number format (date) on column D.
$ws.Rows(19).Insert()

# Populate the newly inserted row 19 with the same data as the (now shifted)
# row 20 except for a new Fecha (date) value, matching the weekly update.
$ws.Range("A19").Value = 10
$ws.Range("B19").Value = "Vega Modelo de Temuco"
$ws.Range("C19").Value = "La Araucanía"
$ws.Range("D19").Value = 44438
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = 300000001
$ws.Range("G19").Value = "Rabanito"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 30
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 6000
$ws.Range("M19").Value = 6000
$ws.Range("N19").Value = "$/docena de paquetes"
$ws.Range("O19").Value = "Provincia de Cautín"
$ws.Range("P19").Value = 500
$ws.Range("Q19").Value = 12
$ws.Range("R19").Value = "Hortaliza"
